$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: translate the Chinese title into English
$ws.Range("A1").Value = "Using the default launch file, run the common_benchmark_node node for 1 hour of data"

# Row 9 ("Gemini 2L") benchmark numbers were updated. These columns store
# numeric-looking values as text in the source data, so a leading
# apostrophe is used to force text entry (matching how the workbook
# already stores them) instead of letting Excel auto-convert to numbers.
$ws.Range("H9").Value = "'108.58"
$ws.Range("I9").Value = "'98.72"
$ws.Range("J9").Value = "'30.06"
$ws.Range("K9").Value = "'"

$ws.Range("P9").Value = "'114.4"
$ws.Range("Q9").Value = "'79.88"
$ws.Range("R9").Value = "'29.73"
$ws.Range("S9").Value = "'155.01"
